$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.352161407470703
$ws.Range("B1").Value = 3.984588861465454
$ws.Range("C1").Value = 3.765624046325684
$ws.Range("D1").Value = 1.683860898017883
$ws.Range("E1").Value = 1.189673066139221
